$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 90909630
$ws.Range("I11").Value = 90909630
$ws.Range("K11").Value = 90909630
$ws.Range("M11").Value = -90909490
$ws.Range("H13").Value = 29998.5
$ws.Range("J13").Value = 29998.5
$ws.Range("L13").Value = 29998.5
$ws.Range("N13").Value = -30336.5
$ws.Range("H17").Value = 31753.455
$ws.Range("J17").Value = 32698.875
$ws.Range("L17").Value = 98096.625
$ws.Range("N17").Value = -98432.625
$ws.Range("H62").Value = 3696.5
$ws.Range("I62").Value = 3120.75
$ws.Range("J62").Value = 5999.5
$ws.Range("K62").Value = 3120.75
$ws.Range("L62").Value = 5999.5
$ws.Range("M62").Value = -2496.75
$ws.Range("N62").Value = -7247.5
$ws.Range("H65").Value = 3696.5
$ws.Range("I65").Value = 3120.75
$ws.Range("J65").Value = 5999.5
$ws.Range("K65").Value = 15603.75
$ws.Range("L65").Value = 29997.5
$ws.Range("M65").Value = -12483.75
$ws.Range("N65").Value = -36237.5
$ws.Range("H80").Value = 2924683.5
$ws.Range("I80").Value = 1961688.8
$ws.Range("J80").Value = 3665448.5
$ws.Range("K80").Value = 5885066.4
$ws.Range("L80").Value = 10996345.5
$ws.Range("M80").Value = -5884068.4
$ws.Range("N80").Value = -10998341.5
$ws.Range("H83").Value = 2924683.5
$ws.Range("I83").Value = 1961688.8
$ws.Range("J83").Value = 3665448.5
$ws.Range("K83").Value = 17655199.2
$ws.Range("L83").Value = 32989036.5
$ws.Range("M83").Value = -17650207.2
$ws.Range("N83").Value = -32999020.5
$ws.Range("H97").Value = 1308.8
$ws.Range("J97").Value = 1398.6666
$ws.Range("L97").Value = 4195.9998
$ws.Range("N97").Value = -5187.9998
$ws.Range("H98").Value = 5953943.5
$ws.Range("I98").Value = 6945518
$ws.Range("K98").Value = 6945518
$ws.Range("M98").Value = -6944020
$ws.Range("H101").Value = 729.0714
$ws.Range("I101").Value = 227.57143
$ws.Range("J101").Value = 1230.5714
$ws.Range("K101").Value = 682.71429
$ws.Range("L101").Value = 3691.7142
$ws.Range("M101").Value = 939.28571
$ws.Range("N101").Value = -6935.7142
$ws.Range("H118").Value = 2935.875
$ws.Range("I118").Value = 3069.7144
$ws.Range("J118").Value = 1999
$ws.Range("K118").Value = 9209.143199999999
$ws.Range("L118").Value = 5997
$ws.Range("M118").Value = -7552.143199999999
$ws.Range("N118").Value = -9311
$ws.Range("H122").Value = 5953943.5
$ws.Range("I122").Value = 6945518
$ws.Range("K122").Value = 20836554
$ws.Range("M122").Value = -20834104
$ws.Range("H125").Value = 1097.1666
$ws.Range("I125").Value = 488.83334
$ws.Range("J125").Value = 1401.3334
$ws.Range("K125").Value = 4399.50006
$ws.Range("L125").Value = 12612.0006
$ws.Range("M125").Value = -1939.50006
$ws.Range("N125").Value = -17532.0006

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9259.514999999999
$ws.Range("I32").Value = 8901.303
$ws.Range("J32").Value = 15170
$ws.Range("K32").Value = 8901.303
$ws.Range("L32").Value = 15170
$ws.Range("M32").Value = -8614.303
$ws.Range("N32").Value = -15744
$ws.Range("H53").Value = 200013940
$ws.Range("J53").Value = 500020000
$ws.Range("L53").Value = 500020000
$ws.Range("N53").Value = -500021364
$ws.Range("H60").Value = 250069730
$ws.Range("I60").Value = 89451
$ws.Range("K60").Value = 89451
$ws.Range("M60").Value = -88718
$ws.Range("H95").Value = 51666.668
$ws.Range("J95").Value = 51666.668
$ws.Range("L95").Value = 51666.668
$ws.Range("N95").Value = -57158.668
$ws.Range("H124").Value = 23199
$ws.Range("J124").Value = 23199
$ws.Range("L124").Value = 23199
$ws.Range("N124").Value = -33019
$ws.Range("H133").Value = 73582.75
$ws.Range("J133").Value = 67110.11
$ws.Range("L133").Value = 67110.11
$ws.Range("N133").Value = -72170.11

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1326.6957
$ws.Range("I86").Value = 986.1111
$ws.Range("J86").Value = 2552.8
$ws.Range("K86").Value = 986.1111
$ws.Range("L86").Value = 2552.8
$ws.Range("M86").Value = 136.8889
$ws.Range("N86").Value = -4798.8
$ws.Range("H88").Value = 32573.25
$ws.Range("J88").Value = 32573.25
$ws.Range("L88").Value = 32573.25
$ws.Range("N88").Value = -33385.25
$ws.Range("H89").Value = 1326.6957
$ws.Range("I89").Value = 986.1111
$ws.Range("J89").Value = 2552.8
$ws.Range("K89").Value = 4930.555499999999
$ws.Range("L89").Value = 12764
$ws.Range("M89").Value = 685.4445000000005
$ws.Range("N89").Value = -23996
$ws.Range("H91").Value = 32573.25
$ws.Range("J91").Value = 32573.25
$ws.Range("L91").Value = 32573.25
$ws.Range("N91").Value = -35381.25
$ws.Range("H105").Value = 539674.4399999999
$ws.Range("I105").Value = 758682.0600000001
$ws.Range("J105").Value = 7798.857
$ws.Range("K105").Value = 758682.0600000001
$ws.Range("L105").Value = 7798.857
$ws.Range("M105").Value = -756935.0600000001
$ws.Range("N105").Value = -11292.857
$ws.Range("H107").Value = 3140.3704
$ws.Range("I107").Value = 3533.5
$ws.Range("J107").Value = 2017.1428
$ws.Range("K107").Value = 3533.5
$ws.Range("L107").Value = 2017.1428
$ws.Range("M107").Value = -1613.5
$ws.Range("N107").Value = -5857.1428

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10408
$ws.Range("I3").Value = 5350.1113
$ws.Range("J3").Value = 19512.2
$ws.Range("K3").Value = 16050.3339
$ws.Range("L3").Value = 58536.60000000001
$ws.Range("M3").Value = -15938.3339
$ws.Range("N3").Value = -58760.60000000001
$ws.Range("H18").Value = 657.75
$ws.Range("I18").Value = 431.8889
$ws.Range("K18").Value = 1295.6667
$ws.Range("M18").Value = -1126.6667
$ws.Range("H32").Value = 3008333.2
$ws.Range("I32").Value = 1000000
$ws.Range("J32").Value = 3677777.8
$ws.Range("K32").Value = 3000000
$ws.Range("L32").Value = 11033333.4
$ws.Range("M32").Value = -2999717
$ws.Range("N32").Value = -11033899.4
$ws.Range("H62").Value = 12247.5
$ws.Range("I62").Value = 9133.4
$ws.Range("J62").Value = 17437.666
$ws.Range("K62").Value = 27400.2
$ws.Range("L62").Value = 52312.99800000001
$ws.Range("M62").Value = -26714.2
$ws.Range("N62").Value = -53684.99800000001
$ws.Range("H63").Value = 22511.143
$ws.Range("I63").Value = 13365
$ws.Range("J63").Value = 27592.334
$ws.Range("K63").Value = 40095
$ws.Range("L63").Value = 82777.00199999999
$ws.Range("M63").Value = -39346
$ws.Range("N63").Value = -84275.00199999999
$ws.Range("H65").Value = 12247.5
$ws.Range("I65").Value = 9133.4
$ws.Range("J65").Value = 17437.666
$ws.Range("K65").Value = 82200.59999999999
$ws.Range("L65").Value = 156938.994
$ws.Range("M65").Value = -78768.59999999999
$ws.Range("N65").Value = -163802.994
$ws.Range("H66").Value = 22511.143
$ws.Range("I66").Value = 13365
$ws.Range("J66").Value = 27592.334
$ws.Range("K66").Value = 120285
$ws.Range("L66").Value = 248331.006
$ws.Range("M66").Value = -116541
$ws.Range("N66").Value = -255819.006
$ws.Range("H94").Value = 19654.125
$ws.Range("I94").Value = 3900
$ws.Range("J94").Value = 21904.715
$ws.Range("K94").Value = 11700
$ws.Range("L94").Value = 65714.145
$ws.Range("M94").Value = -11024
$ws.Range("N94").Value = -67066.145
$ws.Range("H107").Value = 3640984.5
$ws.Range("J107").Value = 4333796
$ws.Range("L107").Value = 13001388
$ws.Range("N107").Value = -13005228
$ws.Range("H109").Value = 13345.75
$ws.Range("I109").Value = 50
$ws.Range("J109").Value = 17777.666
$ws.Range("K109").Value = 150
$ws.Range("L109").Value = 53332.99800000001
$ws.Range("M109").Value = 890
$ws.Range("N109").Value = -55412.99800000001
$ws.Range("H133").Value = 36334.09
$ws.Range("I133").Value = 29292.555
$ws.Range("J133").Value = 68021
$ws.Range("K133").Value = 87877.66500000001
$ws.Range("L133").Value = 204063
$ws.Range("M133").Value = -82817.66500000001
$ws.Range("N133").Value = -214183
$ws.Range("H139").Value = 5120.1562
$ws.Range("I139").Value = 2409.9
$ws.Range("J139").Value = 9637.25
$ws.Range("K139").Value = 7229.700000000001
$ws.Range("L139").Value = 28911.75
$ws.Range("M139").Value = -2089.700000000001
$ws.Range("N139").Value = -39191.75

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H70").Value = 4999.6665
$ws.Range("I70").Value = 4996.5
$ws.Range("J70").Value = 4999.893
$ws.Range("K70").Value = 4996.5
$ws.Range("L70").Value = 4999.893
$ws.Range("M70").Value = -4726.5
$ws.Range("N70").Value = -5539.893
$ws.Range("H73").Value = 4999.6665
$ws.Range("I73").Value = 4996.5
$ws.Range("J73").Value = 4999.893
$ws.Range("K73").Value = 4996.5
$ws.Range("L73").Value = 4999.893
$ws.Range("M73").Value = -4060.5
$ws.Range("N73").Value = -6871.893
$ws.Range("H110").Value = 49998
$ws.Range("J110").Value = 49998
$ws.Range("L110").Value = 49998
$ws.Range("N110").Value = -58178
$ws.Range("H136").Value = 12210.667
$ws.Range("J136").Value = 12210.667
$ws.Range("L136").Value = 36632.001
$ws.Range("N136").Value = -41732.001

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4845.769
$ws.Range("I40").Value = 4246.3
$ws.Range("J40").Value = 6844
$ws.Range("K40").Value = 4246.3
$ws.Range("L40").Value = 6844
$ws.Range("M40").Value = -4110.3
$ws.Range("N40").Value = -7116
$ws.Range("H95").Value = 56667
$ws.Range("J95").Value = 56667
$ws.Range("L95").Value = 56667
$ws.Range("N95").Value = -62159
$ws.Range("H106").Value = 37500
$ws.Range("J106").Value = 37500
$ws.Range("L106").Value = 37500
$ws.Range("N106").Value = -40024

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1569.3572
$ws.Range("I81").Value = 1492.1
$ws.Range("K81").Value = 2984.2
$ws.Range("M81").Value = -1923.2
$ws.Range("H84").Value = 1569.3572
$ws.Range("I84").Value = 1492.1
$ws.Range("K84").Value = 14921
$ws.Range("M84").Value = -9617
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
